$wb = $excel.ActiveWorkbook

# Update "展览" sheet (sheet1) row 2: F2 58 -> 59, G2 29.9 -> 45
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F2").Value = 59
$ws1.Range("G2").Value = 45

# Update "全部类型" sheet (sheet4) row 2: F2 58 -> 59, G2 29.9 -> 45
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F2").Value = 59
$ws4.Range("G2").Value = 45
